$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.388.97"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "2.916.89"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.77"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.04"
$ws.Range("E6").Value = "  -1.79%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.499"
$ws.Range("E8").Value = "  -1.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.01"
$ws.Range("E9").Value = "  +1.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000221"
$ws.Range("E12").Value = "  -1.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.80"
$ws.Range("E13").Value = "  -2.16%  "

$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").Value = "3.398.91"
$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "61.365.21"
$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").Value = "2.915.63"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.58"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "430.70"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.54"
$ws.Range("E20").Value = "  +1.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.665"
$ws.Range("E21").Value = "  -1.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.02"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.07"
$ws.Range("E23").Value = "  -0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  -1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  -3.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.59"
$ws.Range("E26").Value = "  -1.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.15"
$ws.Range("E28").Value = "  -5.42%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.57"
$ws.Range("E29").Value = "  -1.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.84"
$ws.Range("E30").Value = "  -2.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.49"
$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.107"
$ws.Range("E32").Value = "  -0.56%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "0.0₃0872"
$ws.Range("E34").Value = "  +1.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.53"
$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.92"
$ws.Range("E37").Value = "  -3.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("E38").Value = "  -0.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.121"
$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.40"
$ws.Range("E40").Value = "  -1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.44"
$ws.Range("E41").Value = "  +1.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.271"
$ws.Range("E42").Value = "  -5.12%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.677.62"
$ws.Range("E43").Value = "  -0.57%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0339"
$ws.Range("E44").Value = "  -1.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "133.53"
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "351.80"
$ws.Range("E46").Value = "  -6.90%  "

$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.23"
$ws.Range("E48").Value = "  -2.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.104"
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.98"
$ws.Range("E50").Value = "  -1.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.124"
$ws.Range("E51").Value = "  -0.23%  "
